# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — rows 3-6 hold the "想去人数" (interest count) column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 3160
$wsExpo.Range("F4").Value = 47
$wsExpo.Range("F5").Value = 926
$wsExpo.Range("F6").Value = 299

# Sheet "全部类型" (all types) — same underlying rows, but row 7 holds the
# entry that was row 6 on the "展览" sheet because of an extra row in between
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 3160
$wsAll.Range("F4").Value = 47
$wsAll.Range("F5").Value = 926
$wsAll.Range("F7").Value = 299
